$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, shifting existing rows 51..133 down to 52..134
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new data record
$ws.Cells.Item(51, 1).Value = 3
$ws.Cells.Item(51, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(51, 3).Value = "Coquimbo"
$ws.Cells.Item(51, 4).Value = 44495
$ws.Cells.Item(51, 5).Value = 5
$ws.Cells.Item(51, 6).Value = 100112010
$ws.Cells.Item(51, 7).Value = "Achicoria"
$ws.Cells.Item(51, 8).Value = "Sin especificar"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 128
$ws.Cells.Item(51, 11).Value = 6000
$ws.Cells.Item(51, 12).Value = 6500
$ws.Cells.Item(51, 13).Value = 6234
$ws.Cells.Item(51, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(51, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(51, 16).Value = 390
$ws.Cells.Item(51, 17).Value = 16
$ws.Cells.Item(51, 18).Value = "Hortaliza"
